{"js": "// Turn on track-changes so the removal of \"NOT NULL\" is recorded as a\n// tracked deletion (<w:del>/<w:delText>), matching the author's edit.\ncontext.document.changeTrackingMode = \"TrackAll\";\nawait context.sync();\n\n// --- cars table: \"customer_id: INT, NOT NULL FK\" -> \"customer_id: INT,  FK\"\n// (delete just \"NOT NULL\", leaving the surrounding spaces/\"FK\" intact)\nlet outer = context.document.body.search(\"customer_id: INT, NOT NULL FK\", { matchCase: true });\nouter.load(\"items\");\nawait context.sync();\n\nif (outer.items.length > 0) {\n  const inner = outer.items[0].search(\"NOT NULL\", { matchCase: true });\n  inner.load(\"items\");\n  await context.sync();\n  if (inner.items.length > 0) {\n    inner.items[0].delete();\n    await context.sync();\n  }\n}\n\n// --- repair_orders table: \"car_id: INT, NOT NULL, FK\" -> \"car_id: INT, FK\"\n// (delete just \", NOT NULL\", leaving \"car_id: INT\" and \", FK\" intact)\nouter = context.document.body.search(\"car_id: INT, NOT NULL, FK\", { matchCase: true });\nouter.load(\"items\");\nawait context.sync();\n\nif (outer.items.length > 0) {\n  const inner = outer.items[0].search(\", NOT NULL\", { matchCase: true });\n  inner.load(\"items\");\n  await context.sync();\n  if (inner.items.length > 0) {\n    inner.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Attribute the tracked change to the same reviewer as the original edit,\n# and turn on track-changes so the removal of \"NOT NULL\" is recorded as a\n# tracked deletion (<w:del>/<w:delText>).\n$word.UserName = \"H Fillerup\"\n$d = $word.ActiveDocument\n$d.TrackRevisions = $true\n\n# --- cars table: \"customer_id: INT, NOT NULL FK\" -> \"customer_id: INT,  FK\"\n# (delete just \"NOT NULL\", leaving the surrounding spaces/\"FK\" intact)\n$outer = $d.Content\nif ($outer.Find.Execute(\"customer_id: INT, NOT NULL FK\")) {\n    $inner = $outer.Duplicate\n    if ($inner.Find.Execute(\"NOT NULL\")) {\n        $inner.Delete()\n    }\n}\n\n# --- repair_orders table: \"car_id: INT, NOT NULL, FK\" -> \"car_id: INT, FK\"\n# (delete just \", NOT NULL\", leaving \"car_id: INT\" and \", FK\" intact)\n$outer2 = $d.Content\nif ($outer2.Find.Execute(\"car_id: INT, NOT NULL, FK\")) {\n    $inner2 = $outer2.Duplicate\n    if ($inner2.Find.Execute(\", NOT NULL\")) {\n        $inner2.Delete()\n    }\n}\n"}
